# Battery Pack Price: calculate transportation battery pack prices
# endogenously after 2023.
#
# - Calculations sheet: add a new row (4) holding the "endogenous learning"
#   marker values (0) for years 2024-2050 (columns G:AG), highlighted with a
#   light fill, plus a note in row 6 explaining the flag.
# - BPP sheet: years 2024 onward (columns E:AE) now just pull the (zeroed)
#   endogenous value straight from Calculations row 4 instead of computing
#   Calculations row 3 * About!A13 * About!A14.

$wb = $excel.ActiveWorkbook
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsBPP = $wb.Worksheets.Item("BPP")

function Get-ColLetter($col) {
    $letter = ""
    while ($col -gt 0) {
        $rem = ($col - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $col = [int](($col - $rem - 1) / 26)
    }
    return $letter
}

# Highlight fill = Excel's "Orange, Accent 2, Lighter 80%" swatch
# (theme 5 / tint 0.8 in the original file), expressed as its equivalent
# literal RGB (FBE5D6) since this engine normalizes fills to RGB on save.
$highlightColor = 14083579

# --- Calculations sheet: new row 4 (years 2024-2050, columns G:AG) -------
# Mark these years as "use endogenous learning" by writing 0 and flagging
# the cells with a light highlight fill.
$calcRow4 = $wsCalc.Range("G4:AG4")
$calcRow4.Value = 0
$calcRow4.Interior.Color = $highlightColor

# --- Calculations sheet: new row 6 note -----------------------------------
$noteCell = $wsCalc.Range("G6")
$noteCell.Value = "Use endogenous learning in EPS for future years"
$noteCell.Interior.Color = $highlightColor

# --- BPP sheet: years 2024-2050 (columns E:AE) now reference the new ------
# Calculations row 4 directly, with no About multipliers.
for ($col = 5; $col -le 31; $col++) {
    $bppCell = $wsBPP.Cells.Item(2, $col)
    $calcCol = $col + 2
    $calcLetter = Get-ColLetter $calcCol
    $bppCell.Formula = "=Calculations!" + $calcLetter + "4"
}

# --- Cursor/selection bookkeeping (matches the saved workbook view state) -
$wsCalc.Activate() | Out-Null
$wsCalc.Range("I13").Select() | Out-Null

$wsBPP.Activate() | Out-Null
$wsBPP.Range("E2").Select() | Out-Null
